# Fix Typo in week 6/7 PPT
$p = $ppt.ActivePresentation

# --- Slide 2: table-of-contents slide ---------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange

# "실습. 점수창 만들기" -> "실습. 점수 창 만들기"
# Merge the two runs ("점수창" + " 만들기") into a single clean run by
# writing the full replacement text into the second (non-misspelled) run
# and then clearing the first run's text.
$run2b = $tr2.Characters(118, 4)
$run2b.Text = "점수 창 만들기"
$run2a = $tr2.Characters(115, 3)
$run2a.Text = ""

# "World 게임오브젝의 마우스 감지" -> "World 게임오브젝트의 마우스 감지"
$run2c = $tr2.Characters(92, 6)
$run2c.Text = "게임오브젝트의"

# --- Slide 8: Canvas slide ---------------------------------------------
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(2)
$tr8 = $sh8.TextFrame.TextRange

# "하얀선 - Game View의 크기" -> "하얀 테두리 - Game View의 크기"
# Merge "하얀선" + " " into a single clean run.
$run8b = $tr8.Characters(31, 1)
$run8b.Text = "하얀 테두리 "
$run8a = $tr8.Characters(28, 3)
$run8a.Text = ""

# --- Slide 31: scene-move slide -----------------------------------------
$s31 = $p.Slides.Item(31)
$sh31 = $s31.Shapes.Item(2)
$tr31 = $sh31.TextFrame.TextRange

# "Scene = 게임을 시작했을 때 나타나는 Scene" -> "... 불러오는 Scene"
$run31b = $tr31.Characters(92, 16)
$run31b.Text = "게임을 시작했을 때 불러오는 "

# "이 Scene은 빌드용이라고 알려준다" -> "... 설정하는 것"
$run31a = $tr31.Characters(68, 5)
$run31a.Text = " 설정하는 것"
